# Remove waste heat recovery connection (electricity -> CO2 capture).
# The "connections" worksheet had a row describing the "power" chain's
# "waste heat"/"recovered heat" outflow feeding into "CO2 Capture".
# That row is deleted entirely so no waste heat is considered at all;
# the rows below it shift up to fill the gap.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("connections")
$ws.Rows.Item(14).Delete()
